$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 562, shifting existing rows 562:639 down to 563:640
$ws.Rows("562:562").Insert()

# Populate the newly inserted row 562 with the new data record
$ws.Range("A562").Value = 3
$ws.Range("B562").Value = "Femacal de La Calera"
$ws.Range("C562").Value = "Coquimbo"
$ws.Range("D562").Value = 45127
$ws.Range("E562").Value = 5
$ws.Range("F562").Value = 100112031
$ws.Range("G562").Value = "Poroto verde"
$ws.Range("H562").Value = "Magnum"
$ws.Range("I562").Value = "Primera"
$ws.Range("J562").Value = 73
$ws.Range("K562").Value = 26000
$ws.Range("L562").Value = 27000
$ws.Range("M562").Value = 26479
$ws.Range("N562").Value = "$/malla 25 kilos"
$ws.Range("O562").Value = "Región de Arica y Parinacota"
$ws.Range("P562").Value = 1059
$ws.Range("Q562").Value = 25
$ws.Range("R562").Value = "Hortaliza"
